# Fruta / hortaliza, semanal
# The data rows (Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg) got reshuffled between rows. For each affected row we
# capture its original values first (since several rows form permutation cycles),
# then write back the values taken from the row indicated by the mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> source row (value that row `key` should receive comes from row `value`'s
# ORIGINAL contents)
$mapping = @{
    3  = 15
    4  = 22
    5  = 7
    6  = 9
    7  = 3
    8  = 6
    9  = 16
    10 = 20
    11 = 4
    12 = 13
    13 = 21
    14 = 10
    15 = 11
    16 = 5
    20 = 8
    21 = 12
    22 = 23
    23 = 14
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot original values for every row involved before making any changes.
$original = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Apply the new values based on the mapping using the captured originals.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $original[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
